$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove 2 rows of test data (rows 53 and 54) - shifts remaining rows up,
# shrinking the table/used-range from A1:D90 to A1:D88.
$ws.Rows("53:54").Delete()

# The ExternalData_1 defined name is not auto-updated by the row delete
# (it is not linked to the table range dynamically), so fix it manually.
foreach ($n in $wb.Names) {
  if ($n.Name -eq "Sheet1!ExternalData_1") {
    $n.RefersTo = "=Sheet1!`$A`$1:`$C`$88"
  }
}

# Update the selected cell in the sheet view.
$ws.Range("H9").Select()
